$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Daily auto-push: insert a new scraped data row at row 618, shifting the
# existing rows (old 618..659) down by one.
$ws.Rows.Item(618).Insert()

$ws.Cells.Item(618, 1).Value = "'2026/01/13"
$ws.Cells.Item(618, 2).Value = "火"
$ws.Cells.Item(618, 3).Value = 23
$ws.Cells.Item(618, 4).Value = 29
